$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 20:50"

# Countries/provinces data refresh: the source (sorted descending by Casos
# totales) re-ranked several countries, so both the country label in column A
# and the B:H statistics shifted for the affected rows.
$paisesData = @(
    @{ Row = 4; Pais = "Estados Unidos"; Vals = @(266279, 21402, 11983, 247493, 5781, 733, 6803) },
    @{ Row = 24; Pais = "Noruega"; Vals = @(5296, 149, 32, 5206, 96, 8, 58) },
    @{ Row = 41; Pais = "Tailandia"; Vals = @(1978, 103, 612, 1347, 23, 4, 19) },
    @{ Row = 43; Pais = "Grecia"; Vals = @(1613, 69, 78, 1476, 92, 6, 59) },
    @{ Row = 58; Pais = "Egipto"; Vals = @(985, 120, 216, 703, 0, 8, 66) },
    @{ Row = 59; Pais = "Estonia"; Vals = @(961, 103, 48, 901, 16, 1, 12) },
    @{ Row = 60; Pais = "Ucrania"; Vals = @(942, 45, 19, 900, 16, 1, 23) },
    @{ Row = 61; Pais = "Eslovenia"; Vals = @(934, 37, 70, 844, 31, 3, 20) },
    @{ Row = 62; Pais = "Nueva Zelanda"; Vals = @(868, 71, 103, 764, 2, 0, 1) },
    @{ Row = 72; Pais = "Bosnia y Herzegovina"; Vals = @(579, 46, 27, 535, 4, 1, 17) },
    @{ Row = 74; Pais = "Tunez"; Vals = @(495, 40, 5, 472, 30, 4, 18) },
    @{ Row = 77; Pais = "Kazajistan"; Vals = @(464, 29, 29, 429, 6, 3, 6) },
    @{ Row = 81; Pais = "Republica de Macedonia"; Vals = @(430, 46, 20, 398, 8, 1, 12) },
    @{ Row = 98; Pais = "Uzbekistan"; Vals = @(227, 22, 25, 200, 8, 0, 2) },
    @{ Row = 99; Pais = "Honduras"; Vals = @(222, 3, 3, 204, 10, 1, 15) },
    @{ Row = 110; Pais = "Georgia"; Vals = @(155, 21, 27, 128, 6, 0, 0) },
    @{ Row = 117; Pais = "Guadalupe"; Vals = @(130, 2, 24, 99, 14, 1, 7) },
    @{ Row = 126; Pais = "Ruanda"; Vals = @(89, 5, 0, 89, 0, 0, 0) },
    @{ Row = 129; Pais = "Madagascar"; Vals = @(70, 11, 0, 70, 6, 0, 0) },
    @{ Row = 142; Pais = "Polinesia Francesa"; Vals = @(39, 2, 0, 39, 1, 0, 0) },
    @{ Row = 143; Pais = "Zambia"; Vals = @(39, 0, 2, 36, 0, 0, 1) },
    @{ Row = 144; Pais = "Puerto Rico"; Vals = @(39, 0, 1, 36, 0, 0, 2) },
    @{ Row = 145; Pais = "Mali"; Vals = @(39, 3, 0, 36, 0, 0, 3) },
    @{ Row = 160; Pais = "Haiti"; Vals = @(18, 2, 1, 17, 0, 0, 0) },
    @{ Row = 161; Pais = "Nueva Caledonia"; Vals = @(18, 0, 1, 17, 0, 0, 0) },
    @{ Row = 164; Pais = "Benin"; Vals = @(16, 3, 2, 14, 0, 0, 0) },
    @{ Row = 165; Pais = "Siria"; Vals = @(16, 0, 0, 14, 0, 0, 2) },
    @{ Row = 166; Pais = "Mongolia"; Vals = @(14, 0, 2, 12, 0, 0, 0) },
    @{ Row = 167; Pais = "Namibia"; Vals = @(14, 0, 3, 11, 0, 0, 0) },
    @{ Row = 168; Pais = "Santa Lucia"; Vals = @(13, 0, 1, 12, 0, 0, 0) },
    @{ Row = 188; Pais = "Santa Sede"; Vals = @(7, 0, 0, 7, 0, 0, 0) },
    @{ Row = 190; Pais = "Liberia"; Vals = @(7, 1, 0, 7, 0, 0, 0) },
    @{ Row = 191; Pais = "Somalia"; Vals = @(7, 2, 1, 6, 0, 0, 0) },
    @{ Row = 192; Pais = "Cabo Verde"; Vals = @(6, 0, 0, 5, 0, 0, 1) },
    @{ Row = 193; Pais = "San Bartolome"; Vals = @(6, 0, 1, 5, 0, 0, 0) },
    @{ Row = 194; Pais = "Nepal"; Vals = @(6, 0, 1, 5, 0, 0, 0) },
    @{ Row = 195; Pais = "Mauritania"; Vals = @(6, 0, 2, 3, 0, 0, 1) },
    @{ Row = 196; Pais = "Islas Turcas y Caicos"; Vals = @(5, 0, 0, 5, 0, 0, 0) },
    @{ Row = 202; Pais = "Burundi"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 203; Pais = "Malaui"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 204; Pais = "Anguila"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 205; Pais = "Islas Virgenes Britanicas"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 209; Pais = "Timor Oriental"; Vals = @(1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 210; Pais = "Papua Nueva Guinea"; Vals = @(1, 0, 0, 1, 0, 0, 0) }
)

foreach ($entry in $paisesData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Pais
    $v = $entry.Vals
    for ($i = 0; $i -lt $v.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $v[$i]
    }
}
